# Weekly update: insert the latest week's data (date 2022-03-21 = serial 44641)
# for "Cultivar IV Región" just above the previous week's block (which started
# at row 154), shifting all subsequent rows down by 4. This mirrors Excel's
# normal "insert rows, enter new data" workflow for a fresh week of prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 154 (rows 154-210 shift down to 158-214).
$ws.Rows("154:157").Insert()

# --- New row 154: Especial ---
$ws.Range("A154").Value2 = 6
$ws.Range("B154").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C154").Value2 = "Metropolitana"
$ws.Range("D154").Value2 = 44641
$ws.Range("E154").Value2 = 13
$ws.Range("F154").Value2 = 100112043
$ws.Range("G154").Value2 = "Pepino dulce"
$ws.Range("H154").Value2 = "Cultivar IV Región"
$ws.Range("I154").Value2 = "Especial"
$ws.Range("J154").Value2 = 150
$ws.Range("K154").Value2 = 14000
$ws.Range("L154").Value2 = 14000
$ws.Range("M154").Value2 = 14000
$ws.Range("N154").Value2 = "$/bandeja 18 kilos"
$ws.Range("O154").Value2 = "Provincia de Limarí"
$ws.Range("P154").Value2 = 778
$ws.Range("Q154").Value2 = 18
$ws.Range("R154").Value2 = "Hortaliza"

# --- New row 155: Primera ---
$ws.Range("A155").Value2 = 6
$ws.Range("B155").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C155").Value2 = "Metropolitana"
$ws.Range("D155").Value2 = 44641
$ws.Range("E155").Value2 = 13
$ws.Range("F155").Value2 = 100112043
$ws.Range("G155").Value2 = "Pepino dulce"
$ws.Range("H155").Value2 = "Cultivar IV Región"
$ws.Range("I155").Value2 = "Primera"
$ws.Range("J155").Value2 = 260
$ws.Range("K155").Value2 = 13000
$ws.Range("L155").Value2 = 13000
$ws.Range("M155").Value2 = 13000
$ws.Range("N155").Value2 = "$/bandeja 18 kilos"
$ws.Range("O155").Value2 = "Provincia de Limarí"
$ws.Range("P155").Value2 = 722
$ws.Range("Q155").Value2 = 18
$ws.Range("R155").Value2 = "Hortaliza"

# --- New row 156: Segunda ---
$ws.Range("A156").Value2 = 6
$ws.Range("B156").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C156").Value2 = "Metropolitana"
$ws.Range("D156").Value2 = 44641
$ws.Range("E156").Value2 = 13
$ws.Range("F156").Value2 = 100112043
$ws.Range("G156").Value2 = "Pepino dulce"
$ws.Range("H156").Value2 = "Cultivar IV Región"
$ws.Range("I156").Value2 = "Segunda"
$ws.Range("J156").Value2 = 130
$ws.Range("K156").Value2 = 10000
$ws.Range("L156").Value2 = 10000
$ws.Range("M156").Value2 = 10000
$ws.Range("N156").Value2 = "$/bandeja 18 kilos"
$ws.Range("O156").Value2 = "Provincia de Limarí"
$ws.Range("P156").Value2 = 556
$ws.Range("Q156").Value2 = 18
$ws.Range("R156").Value2 = "Hortaliza"

# --- New row 157: Tercera ---
$ws.Range("A157").Value2 = 6
$ws.Range("B157").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C157").Value2 = "Metropolitana"
$ws.Range("D157").Value2 = 44641
$ws.Range("E157").Value2 = 13
$ws.Range("F157").Value2 = 100112043
$ws.Range("G157").Value2 = "Pepino dulce"
$ws.Range("H157").Value2 = "Cultivar IV Región"
$ws.Range("I157").Value2 = "Tercera"
$ws.Range("J157").Value2 = 90
$ws.Range("K157").Value2 = 8000
$ws.Range("L157").Value2 = 8000
$ws.Range("M157").Value2 = 8000
$ws.Range("N157").Value2 = "$/bandeja 18 kilos"
$ws.Range("O157").Value2 = "Provincia de Limarí"
$ws.Range("P157").Value2 = 444
$ws.Range("Q157").Value2 = 18
$ws.Range("R157").Value2 = "Hortaliza"
